# Update Leve profit-calc columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job sheets, per refreshed market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 18
$ws.Range("H18").Value = 1420
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# row 28
$ws.Range("H28").Value = 5171.1816
$ws.Range("I28").Value = 1823
$ws.Range("J28").Value = 7961.3335
$ws.Range("K28").Value = 1823
$ws.Range("L28").Value = 7961.3335
$ws.Range("M28").Value = -1338
$ws.Range("N28").Value = -8931.333500000001

# row 62
$ws.Range("H62").Value = 11066.667
$ws.Range("J62").Value = 11066.667
$ws.Range("L62").Value = 11066.667
$ws.Range("N62").Value = -12314.667

# row 65
$ws.Range("H65").Value = 11066.667
$ws.Range("J65").Value = 11066.667
$ws.Range("L65").Value = 55333.335
$ws.Range("N65").Value = -61573.335

# row 112
$ws.Range("H112").Value = 1905
$ws.Range("I112").Value = 1695
$ws.Range("J112").Value = 1981.3636
$ws.Range("K112").Value = 5085
$ws.Range("L112").Value = 5944.0908
$ws.Range("M112").Value = -3977
$ws.Range("N112").Value = -8160.0908

# row 132
$ws.Range("H132").Value = 12040.083
$ws.Range("I132").Value = 11159.35
$ws.Range("J132").Value = 16443.75
$ws.Range("K132").Value = 33478.05
$ws.Range("L132").Value = 49331.25
$ws.Range("M132").Value = -30948.05
$ws.Range("N132").Value = -54391.25

# row 135
$ws.Range("H135").Value = 873.6667
$ws.Range("I135").Value = 895.375
$ws.Range("K135").Value = 8058.375
$ws.Range("M135").Value = -5523.375

# row 137
$ws.Range("H137").Value = 1287.2727
$ws.Range("I137").Value = 997.1429000000001
$ws.Range("J137").Value = 1795
$ws.Range("K137").Value = 2991.4287
$ws.Range("L137").Value = 5385
$ws.Range("M137").Value = -441.4287000000004
$ws.Range("N137").Value = -10485

# row 138
$ws.Range("H138").Value = 2906.7827
$ws.Range("I138").Value = 743.9091
$ws.Range("J138").Value = 4889.4165
$ws.Range("K138").Value = 2231.7273
$ws.Range("L138").Value = 14668.2495
$ws.Range("M138").Value = 2908.2727
$ws.Range("N138").Value = -24948.2495

$ws = $wb.Worksheets.Item("ARM")
# row 39
$ws.Range("H39").Value = 4218.5
$ws.Range("I39").Value = 3062.2
$ws.Range("K39").Value = 3062.2
$ws.Range("M39").Value = -2542.2

# row 74
$ws.Range("H74").Value = 2826.0908
$ws.Range("I74").Value = 2608.7
$ws.Range("K74").Value = 2608.7
$ws.Range("M74").Value = -1734.7

# row 77
$ws.Range("H77").Value = 2826.0908
$ws.Range("I77").Value = 2608.7
$ws.Range("K77").Value = 13043.5
$ws.Range("M77").Value = -8675.5

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 2838.4285
$ws.Range("J134").Value = 4400
$ws.Range("L134").Value = 13200
$ws.Range("N134").Value = -18270

$ws = $wb.Worksheets.Item("CRP")
# row 62
$ws.Range("H62").Value = 8601
$ws.Range("I62").Value = 8601
$ws.Range("K62").Value = 8601
$ws.Range("M62").Value = -7977

# row 65
$ws.Range("H65").Value = 8601
$ws.Range("I65").Value = 8601
$ws.Range("K65").Value = 43005
$ws.Range("M65").Value = -39885

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1277.3572
$ws.Range("I5").Value = 1004.25
$ws.Range("J5").Value = 1641.5
$ws.Range("K5").Value = 3012.75
$ws.Range("L5").Value = 4924.5
$ws.Range("M5").Value = -2900.75
$ws.Range("N5").Value = -5148.5

# row 68
$ws.Range("H68").Value = 1862.25
$ws.Range("I68").Value = 1449
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 4347
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -3536
$ws.Range("N68").Value = -7622

# row 71
$ws.Range("H71").Value = 1862.25
$ws.Range("I71").Value = 1449
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13041
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = -8985
$ws.Range("N71").Value = -26112

# row 117
$ws.Range("H117").Value = 171.4
$ws.Range("J117").Value = 166.71428
$ws.Range("L117").Value = 500.14284
$ws.Range("N117").Value = -7384.14284

# row 119
$ws.Range("H119").Value = 5676.3335

# row 131
$ws.Range("H131").Value = 1966.2727
$ws.Range("I131").Value = 1307.25
$ws.Range("K131").Value = 3921.75
$ws.Range("M131").Value = 1118.25

# row 135
$ws.Range("H135").Value = 1277.3572
$ws.Range("I135").Value = 1004.25
$ws.Range("J135").Value = 1641.5
$ws.Range("K135").Value = 9038.25
$ws.Range("L135").Value = 14773.5
$ws.Range("M135").Value = -6503.25
$ws.Range("N135").Value = -19843.5

# row 140
$ws.Range("H140").Value = 2406.75
$ws.Range("I140").Value = 1542.3334
$ws.Range("K140").Value = 4627.0002
$ws.Range("M140").Value = 552.9997999999996

$ws = $wb.Worksheets.Item("GSM")
# row 95
$ws.Range("H95").Value = 24499.6
$ws.Range("J95").Value = 24499.6
$ws.Range("L95").Value = 24499.6
$ws.Range("N95").Value = -29991.6

# row 107
$ws.Range("H107").Value = 1360.7858
$ws.Range("I107").Value = 1398.1
$ws.Range("J107").Value = 1267.5
$ws.Range("K107").Value = 1398.1
$ws.Range("L107").Value = 1267.5
$ws.Range("M107").Value = 521.9000000000001
$ws.Range("N107").Value = -5107.5

# row 132
$ws.Range("H132").Value = 8997.333000000001
$ws.Range("I132").Value = 7012
$ws.Range("K132").Value = 21036
$ws.Range("M132").Value = -18506

$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 5249.75
$ws.Range("I68").Value = 2399.6
$ws.Range("K68").Value = 2399.6
$ws.Range("M68").Value = -1650.6

# row 71
$ws.Range("H71").Value = 5249.75
$ws.Range("I71").Value = 2399.6
$ws.Range("K71").Value = 11998
$ws.Range("M71").Value = -8254

# row 93
$ws.Range("H93").Value = 524.25
$ws.Range("I93").Value = 524.25
$ws.Range("K93").Value = 524.25
$ws.Range("M93").Value = 723.75

$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 260.5
$ws.Range("I107").Value = 262
$ws.Range("J107").Value = 250
$ws.Range("K107").Value = 786
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1134
$ws.Range("N107").Value = -4590

# row 126
$ws.Range("H126").Value = 4573.0415
$ws.Range("J126").Value = 6916.6665
$ws.Range("L126").Value = 20749.9995
$ws.Range("N126").Value = -25689.9995

# row 136
$ws.Range("H136").Value = 4226.2
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 5299.125
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 15897.375
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -20997.375
